$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.221774935722351
$ws.Range("B1").Value = 2.723650217056274
$ws.Range("C1").Value = 4.478193759918213
$ws.Range("D1").Value = 2.126012325286865
$ws.Range("E1").Value = 1.160516858100891
